# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# - Updates "Valor Mora" total (E11) and "Cant. Periodos" (F13) to reflect
#   the newly added period.
# - Inserts a new table row for period "2508" (row 21), pushing the closing
#   signature block down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update summary figures -------------------------------------------------
# Valor Mora total goes from 5 periods (284700) to 6 periods (341640)
$ws.Range("E11").Value = 341640
# Cant. Periodos goes from 5 to 6
$ws.Range("F13").Value = 6

# --- Insert the new period row ----------------------------------------------
# Insert a new row right after the last existing data row (20); this pushes
# the trailing signature block (old rows 25-26) down to rows 26-27.
$ws.Rows.Item(21).Insert()

# The new row should carry the "closing" (bottom-border) row style that row
# 20 currently has, since it is now the last row of the table.
$ws.Range("B20:J20").Copy()
$ws.Range("B21:J21").PasteSpecial(-4122)   # xlPasteFormats

# Row 20 is no longer the last row, so it reverts to the regular/middle
# row style, matching rows 16-19.
$ws.Range("B19:J19").Copy()
$ws.Range("B20:J20").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0

# Fill in the new row's values for period 2508
$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1047398728"
$ws.Range("D21").Value = "JESICA PAOLA RODELO JIMENEZ"
$ws.Range("E21").Value = "2508"
$ws.Range("F21").Value = 56940
$ws.Range("G21").Value = 1423500
